$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 5 entirely
$ws.Range("A5:AU5").EntireRow.Delete() | Out-Null

# Delete columns AF:AU entirely
$ws.Range("AF1:AU4").EntireColumn.Delete() | Out-Null

# Update row 3 values V3:AE3
$ws.Range("V3").Value = 49
$ws.Range("W3").Value = 83
$ws.Range("X3").Value = 90
$ws.Range("Y3").Value = 91
$ws.Range("Z3").Value = 92
$ws.Range("AA3").Value = 93
$ws.Range("AB3").Value = 94
$ws.Range("AC3").Value = 95
$ws.Range("AD3").Value = 96
$ws.Range("AE3").Value = 97
